# Apply the "Add files via upload" edit to slide 5 ("BREVE HISTORIA"):
#  - Bold the title text
#  - Resize/reposition the content placeholder and insert a new
#    "BitKeeper" line right after the first blank paragraph
#  - Reposition the four screenshot pictures

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Title: make the run bold ---------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Font.Bold = 1

# --- Content placeholder: move/resize + add "BitKeeper" line --------
$content = $s.Shapes.Item(2)
$content.Left = 45.763149606299216
$content.Top = 55.94582677165354
$content.Width = 868.4736220472441
$content.Height = 376.054094488189

$tr = $content.TextFrame.TextRange
# Paragraph 2 is the blank line right after "Copias añadiendo fechas";
# paragraph 3 is "CVS". Insert the new line between them.
$cvsPara = $tr.Paragraphs(3, 1)
$cvsPara.InsertBefore("BitKeeper`r")

# --- Pictures: reposition (sizes stay the same) ----------------------
$pic1 = $s.Shapes.Item(3)
$pic1.Left = 363.53267716535436
$pic1.Top = 57.93267716535433

$pic2 = $s.Shapes.Item(4)
$pic2.Left = 314.03267716535436
$pic2.Top = 168.59826771653545

$pic3 = $s.Shapes.Item(5)
$pic3.Left = 242.4076377952756
$pic3.Top = 242.68960629921258

$pic4 = $s.Shapes.Item(6)
$pic4.Left = 326.22937007874015
$pic4.Top = 340.0677952755905
